# Updated symbol list on Mon Dec 26 15:41:56 UTC 2022 with GitHub Actions
# Refresh prices/volume stats for existing coins and shift coin listings
# in rows 10-18 up one rank (row 10's old "One" entry moves to row 18,
# with WazirX/MandalaExchangeToken/... each shifting up into the row above).
# All Price-column values are written with a leading apostrophe so they
# stay stored as text (matching the original inlineStr cells) instead of
# being auto-coerced to numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.83"
$ws.Range("D3").Value = "'23.05"
$ws.Range("D4").Value = "'5.428"
$ws.Range("D5").Value = "'0.05891"
$ws.Range("D6").Value = "'3.446"
$ws.Range("D7").Value = "'6.543"
$ws.Range("D8").Value = "'0.8118"
$ws.Range("D9").Value = "'0.9597"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1417"
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.07449"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = "'0.03257"
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.03063"
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = "'0.09336"
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = "'3.868"
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = "'0.001568"
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = "'0.04680"
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = "'0.0005948"
$ws.Range("E18").Value = '17OneONE'
$ws.Range("D19").Value = "'0.005857"
$ws.Range("D20").Value = "'0.001252"
$ws.Range("D24").Value = "'2.129"
$ws.Range("D25").Value = "'0.3229"
$ws.Range("D27").Value = "'0.0002287"
$ws.Range("D40").Value = "'0.03929"
$ws.Range("D41").Value = "'0.006182"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D43").Value = "'0.002544"
$ws.Range("D44").Value = "'0.009154"
$ws.Range("E44").Value = '43LocalTradersLCTBestin24h'
$ws.Range("D45").Value = "'0.00005199"
$ws.Range("D47").Value = "'0.7310"
$ws.Range("D48").Value = "'0.002297"
